# Apply odds/score updates from the 2024-10-17 FlashScore refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("AI2").Value = 26
$ws.Range("G2").Value = 1.73
$ws.Range("I2").Value = 5.25
$ws.Range("U2").Value = 1.95
$ws.Range("V2").Value = 1.8
$ws.Range("W2").Value = 6.5

# Row 5
$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 11
$ws.Range("O5").Value = 1.29
$ws.Range("P5").Value = 3.5
$ws.Range("Q5").Value = 1.95
$ws.Range("R5").Value = 1.9

# Row 7
$ws.Range("AA7").Value = 21
$ws.Range("AB7").Value = 25
$ws.Range("AC7").Value = 10.25
$ws.Range("AE7").Value = 11.5
$ws.Range("AH7").Value = 9.25
$ws.Range("AI7").Value = 14
$ws.Range("AK7").Value = 30
$ws.Range("AM7").Value = 26
$ws.Range("AN7").Value = 4.7
$ws.Range("AT7").Value = 2.75
$ws.Range("AW7").Value = 4.6
$ws.Range("G7").Value = 2.65
$ws.Range("H7").Value = 3.1
$ws.Range("I7").Value = 2.57
$ws.Range("J7").Value = 3.2
$ws.Range("L7").Value = 3.15
$ws.Range("M7").Value = 9.800000000000001
$ws.Range("N7").Value = 1.04
$ws.Range("O7").Value = 1.25
$ws.Range("P7").Value = 3.2
$ws.Range("Q7").Value = 1.8
$ws.Range("R7").Value = 1.91
$ws.Range("S7").Value = 1.35
$ws.Range("T7").Value = 3.04
$ws.Range("U7").Value = 1.57
$ws.Range("V7").Value = 2.12

# Row 8
$ws.Range("AC8").Value = 11
$ws.Range("AG8").Value = 201
$ws.Range("AH8").Value = 12
$ws.Range("AO8").Value = 11
$ws.Range("AP8").Value = 21
$ws.Range("AW8").Value = 5.5
$ws.Range("AX8").Value = 19
$ws.Range("BC8").Value = 126
$ws.Range("G8").Value = 2
$ws.Range("H8").Value = 3.25
$ws.Range("I8").Value = 3.75
$ws.Range("J8").Value = 2.63
$ws.Range("K8").Value = 2.2
$ws.Range("M8").Value = 1.05
$ws.Range("N8").Value = 11
$ws.Range("O8").Value = 1.25
$ws.Range("P8").Value = 3.75
$ws.Range("Q8").Value = 1.88
$ws.Range("R8").Value = 1.98
$ws.Range("W8").Value = 8
$ws.Range("Y8").Value = 9

# Row 9
$ws.Range("AA9").Value = 12
$ws.Range("AH9").Value = 17
$ws.Range("AJ9").Value = 19
$ws.Range("AP9").Value = 17
$ws.Range("BC9").Value = 126
$ws.Range("G9").Value = 1.55
$ws.Range("H9").Value = 3.9
$ws.Range("I9").Value = 6
$ws.Range("K9").Value = 2.38
$ws.Range("M9").Value = 1.03
$ws.Range("N9").Value = 15
$ws.Range("O9").Value = 1.2
$ws.Range("P9").Value = 4.33
$ws.Range("Q9").Value = 1.67
$ws.Range("R9").Value = 2.15
$ws.Range("U9").Value = 1.73
$ws.Range("V9").Value = 2
$ws.Range("W9").Value = 8

# Row 11
$ws.Range("N11").Value = 8

# Row 12
$ws.Range("AI12").Value = 15
$ws.Range("AO12").Value = 13
$ws.Range("AP12").Value = 26
$ws.Range("AQ12").Value = 41
$ws.Range("AX12").Value = 21
$ws.Range("G12").Value = 2.25
$ws.Range("I12").Value = 3.2
$ws.Range("J12").Value = 3.1
$ws.Range("N12").Value = 7.5
$ws.Range("X12").Value = 9.5

# Row 13
$ws.Range("AB13").Value = 34
$ws.Range("AC13").Value = 9.5
$ws.Range("AD13").Value = 8.5
$ws.Range("AH13").Value = 21
$ws.Range("AJ13").Value = 29
$ws.Range("AM13").Value = 67
$ws.Range("AN13").Value = 3.2
$ws.Range("AO13").Value = 6.5
$ws.Range("AQ13").Value = 19
$ws.Range("AU13").Value = 10
$ws.Range("AW13").Value = 9
$ws.Range("AZ13").Value = 201
$ws.Range("G13").Value = 1.38
$ws.Range("I13").Value = 9.5
$ws.Range("L13").Value = 8
$ws.Range("M13").Value = 1.06
$ws.Range("N13").Value = 9.5
$ws.Range("O13").Value = 1.25
$ws.Range("P13").Value = 3.75
$ws.Range("Q13").Value = 1.9
$ws.Range("R13").Value = 1.95
$ws.Range("U13").Value = 2.2
$ws.Range("V13").Value = 1.62
$ws.Range("W13").Value = 6
$ws.Range("X13").Value = 6
$ws.Range("Z13").Value = 8.5
